# Centers sheet gains several new columns: Center Name, Center Num,
# Facility Address Line 1, Facility City Name (x2), Facility Postal Code,
# then the existing Deliv Center Capac / Latitude / Longitude columns move
# right to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Centers")

# ---- Header row ----
# B1 ("Facility Name") is already correct and stays untouched.
$ws.Cells.Item(1, 3).Value = "Center Name"
$ws.Cells.Item(1, 4).Value = "Center Num"
$ws.Cells.Item(1, 5).Value = "Facility Address Line 1"
$ws.Cells.Item(1, 6).Value = "Facility City Name"
$ws.Cells.Item(1, 7).Value = "Facility City Name"
$ws.Cells.Item(1, 8).Value = "Facility Postal Code"
$ws.Cells.Item(1, 9).Value = "Deliv Center Capac"
$ws.Cells.Item(1, 10).Value = "Latitude"
$ws.Cells.Item(1, 11).Value = "Longitude"

# Copy header style (bold / bordered / centered) from the pre-existing
# header cells onto the newly added header cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Data rows ----
# Row 2 : MYKAWA
$ws.Cells.Item(2, 1).Value = 13
$ws.Cells.Item(2, 2).Value = "MYKAWA"
$ws.Cells.Item(2, 3).Value = "MYKAWA EAM"
$ws.Cells.Item(2, 4).Value = 7723
$ws.Cells.Item(2, 5).Value = "7110 MYKAWA ROAD"
$ws.Cells.Item(2, 6).Value = "HOUSTON"
$ws.Cells.Item(2, 7).Value = "HOUSTON"
$ws.Cells.Item(2, 8).Value = 77033
$ws.Cells.Item(2, 9).Value = 999999999
$ws.Cells.Item(2, 10).Value = 29.67578534220857
$ws.Cells.Item(2, 11).Value = -95.32125610590822

# Row 3 : STAFFORD
$ws.Cells.Item(3, 1).Value = 19
$ws.Cells.Item(3, 2).Value = "STAFFORD"
$ws.Cells.Item(3, 3).Value = "STAFFORD-KATY"
$ws.Cells.Item(3, 4).Value = 7741
$ws.Cells.Item(3, 5).Value = "13922 STAFFORD ROAD"
$ws.Cells.Item(3, 6).Value = "STAFFORD"
$ws.Cells.Item(3, 7).Value = "STAFFORD"
$ws.Cells.Item(3, 8).Value = 77477
$ws.Cells.Item(3, 9).Value = 999999999
$ws.Cells.Item(3, 10).Value = 29.61935353306665
$ws.Cells.Item(3, 11).Value = -95.55583702616946

# Row 4 : HOUSTON (CANINO) HUB
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "HOUSTON (CANINO) HUB"
$ws.Cells.Item(4, 3).Value = "SWEETWATER-CENTRAL"
$ws.Cells.Item(4, 4).Value = 7707
$ws.Cells.Item(4, 5).Value = "8330 SWEETWATER LANE"
$ws.Cells.Item(4, 6).Value = "HOUSTON"
$ws.Cells.Item(4, 7).Value = "HOUSTON"
$ws.Cells.Item(4, 8).Value = 77037
$ws.Cells.Item(4, 9).Value = 999999999
$ws.Cells.Item(4, 10).Value = 29.88062033368866
$ws.Cells.Item(4, 11).Value = -95.40848234626999
